# v3.0 update FCI 27/1/2023
# Adds a new date column (C) with updated figures and re-sorts the fund
# rows alphabetically, pushing the "avg"/"total" summary rows to the
# bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column header (copy formatting from B1, then set its value) ---
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("C1").Value = "13-01-2023"

# --- Final layout for rows 2-10: label, column-B (06-01-2023) value, column-C (13-01-2023) value ---
$rows = @(
    @("Alpha Latam",            469.33,               501.48),
    @("Delta Acciones",         0,                    0),
    @("Fima Acciones",          3528.45,              4228.83),
    @("Fima PB Acciones",       11804.95,             11113.28),
    @("HF Acciones Argentinas", 20195.57,             19958.94),
    @("MAF",                    0,                    0),
    @("Supefondo RV",           21964.31,             27091.71),
    @("avg",                    8280.370000000001,    8984.889999999999),
    @("total",                  57962.61,             62894.24)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}
